$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.286.62"
$ws.Range("E2").Value = "  -2.37%  "
$ws.Range("D3").Value = "3.820.76"
$ws.Range("E3").Value = "  -2.46%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.54"
$ws.Range("E5").Value = "  -1.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.41"
$ws.Range("E6").Value = "  -0.36%  "
$ws.Range("D7").Value = "3.822.94"
$ws.Range("E7").Value = "  -2.39%  "
$ws.Range("E9").Value = "  -1.77%  "
$ws.Range("E10").Value = "  -3.03%  "
$ws.Range("E11").Value = "  +1.06%  "
$ws.Range("E12").Value = "  -2.40%  "
$ws.Range("E13").Value = "  +1.73%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.01"
$ws.Range("E14").Value = "  -3.53%  "
$ws.Range("D15").Value = "4.452.33"
$ws.Range("E15").Value = "  -2.71%  "
$ws.Range("D16").Value = "3.804.14"
$ws.Range("E16").Value = "  -2.76%  "
$ws.Range("D17").Value = "68.237.99"
$ws.Range("E17").Value = "  -2.40%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.53"
$ws.Range("E18").Value = "  -1.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.41"
$ws.Range("E19").Value = "  -2.82%  "
$ws.Range("E20").Value = "  -0.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.10"
$ws.Range("E21").Value = "  -0.66%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "467.66"
$ws.Range("E22").Value = "  -5.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.734"
$ws.Range("E23").Value = "  -1.73%  "
$ws.Range("E24").Value = "  -4.88%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.91"
$ws.Range("E25").Value = "  -3.28%  "
$ws.Range("E26").Value = "  -2.45%  "
$ws.Range("E27").Value = "  -1.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.02"
$ws.Range("E28").Value = "  -1.07%  "
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("E30").Value = "  -1.21%  "
$ws.Range("D31").Value = "3.968.63"
$ws.Range("E31").Value = "  -2.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.66"
$ws.Range("E32").Value = "  -2.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.43"
$ws.Range("E33").Value = "  -2.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.31"
$ws.Range("E34").Value = "  -5.14%  "
$ws.Range("E35").Value = "  -0.91%  "
$ws.Range("D36").Value = "3.780.64"
$ws.Range("E36").Value = "  -2.51%  "
$ws.Range("E37").Value = "  -3.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.69"
$ws.Range("E38").Value = "  +12.28%  "
$ws.Range("E39").Value = "  -1.20%  "
$ws.Range("E40").Value = "  -2.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.93"
$ws.Range("E41").Value = "  -3.42%  "
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.315"
$ws.Range("E43").Value = "  -4.45%  "
$ws.Range("E44").Value = "  -6.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.73"
$ws.Range("E45").Value = "  +0.82%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.000296"
$ws.Range("E46").Value = "  +8.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "418.45"
$ws.Range("E47").Value = "  -4.00%  "
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "46.91"
$ws.Range("E49").Value = "  -2.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "26.29"
$ws.Range("E50").Value = "  +3.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "141.46"
$ws.Range("E51").Value = "  -1.06%  "
